$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.875.33"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "3.721.04"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'619.86"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").Value = "'181.34"
$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("D7").Value = "3.721.01"
$ws.Range("E7").Value = "  -2.09%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -3.38%  "

$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -1.34%  "

$ws.Range("D12").Value = "'0.482"
$ws.Range("E12").Value = "  -4.90%  "

$ws.Range("D13").Value = "'40.17"
$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").Value = "'0.0000255"
$ws.Range("E14").Value = "  -2.91%  "

$ws.Range("D15").Value = "4.323.32"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").Value = "3.698.95"
$ws.Range("E16").Value = "  -2.70%  "

$ws.Range("D17").Value = "69.818.04"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").Value = "'16.40"
$ws.Range("E20").Value = "  -3.45%  "

$ws.Range("D21").Value = "'501.70"
$ws.Range("E21").Value = "  -4.47%  "

$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").Value = "'0.726"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").Value = "'2.53"
$ws.Range("E24").Value = "  +1.52%  "

$ws.Range("D25").Value = "'86.37"
$ws.Range("E25").Value = "  -2.01%  "

$ws.Range("D26").Value = "'13.04"
$ws.Range("E26").Value = "  -4.50%  "

$ws.Range("D27").Value = "'11.14"
$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("D28").Value = "'0.0000132"
$ws.Range("E28").Value = "  +6.35%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'2.46"
$ws.Range("E30").Value = "  -1.96%  "

$ws.Range("D31").Value = "'2.92"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("D32").Value = "'7.97"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").Value = "'30.33"
$ws.Range("E33").Value = "  -6.81%  "

$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").Value = "'6.08"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").Value = "'0.138"
$ws.Range("E38").Value = "  +4.51%  "

$ws.Range("D39").Value = "'0.342"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.03"
$ws.Range("E40").Value = "  +7.66%  "

$ws.Range("D41").Value = "'2.07"
$ws.Range("E41").Value = "  -5.59%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'50.00"
$ws.Range("E42").Value = "  -2.73%  "

$ws.Range("D43").Value = "'429.37"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").Value = "'44.16"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Value = "'8.61"
$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("D46").Value = "2.950.88"
$ws.Range("E46").Value = "  -6.57%  "

$ws.Range("D47").Value = "'0.0361"
$ws.Range("E47").Value = "  -2.40%  "

$ws.Range("D48").Value = "'27.42"
$ws.Range("E48").Value = "  -1.47%  "

$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").Value = "'136.82"
$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  -1.85%  "
